$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Air_Sabre_OneWay")
$ws1.Rows("3:5").Delete()
$ws2 = $wb.Worksheets.Item("Air_Sabre_RoundTrip")
$ws2.Rows("3:5").Delete()

$ws1.Range("B2").Value = "Search|AddToCart|CHECKOUTTRIP|LOGIN|ENTERPAXINFO|CONFIRMPAXINFO"
$ws2.Range("B2").Value = "Search|AddToCart|CHECKOUTTRIP|LOGIN|ENTERPAXINFO|CONFIRMPAXINFO"
Write-Host "Done"
